$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Hoja1" to "Words"
$ws.Name = "Words"

# Append new vocabulary rows (Kleuren/Colors + Dieren/Animals categories).
# Cell writes are ordered to match the original authoring sequence so that
# shared-string indices line up with the source workbook.
$ws.Range("A102").Value = "Kleuren"
$ws.Range("B102").Value = 1
$ws.Range("C102").Value = "rood"
$ws.Range("D102").Value = "rojo"
$ws.Range("A103").Value = "Kleuren"
$ws.Range("B103").Value = 1
$ws.Range("C103").Value = "oranje"
$ws.Range("D103").Value = "naranja"
$ws.Range("A104").Value = "Kleuren"
$ws.Range("B104").Value = 1
$ws.Range("C104").Value = "geel"
$ws.Range("D104").Value = "amarillo"
$ws.Range("A105").Value = "Kleuren"
$ws.Range("B105").Value = 1
$ws.Range("C105").Value = "groen"
$ws.Range("D105").Value = "verde"
$ws.Range("A106").Value = "Kleuren"
$ws.Range("B106").Value = 1
$ws.Range("C106").Value = "blauw"
$ws.Range("D106").Value = "azul"
$ws.Range("A107").Value = "Kleuren"
$ws.Range("B107").Value = 1
$ws.Range("C107").Value = "paars"
$ws.Range("A108").Value = "Kleuren"
$ws.Range("B108").Value = 1
$ws.Range("C108").Value = "roze"
$ws.Range("D108").Value = "rosa"
$ws.Range("D107").Value = "morado"
$ws.Range("A109").Value = "Kleuren"
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "bruin"
$ws.Range("D109").Value = "café"
$ws.Range("A110").Value = "Kleuren"
$ws.Range("B110").Value = 1
$ws.Range("C110").Value = "zwart"
$ws.Range("D110").Value = "negro"
$ws.Range("A111").Value = "Kleuren"
$ws.Range("B111").Value = 1
$ws.Range("C111").Value = "wit"
$ws.Range("D111").Value = "blanco"
$ws.Range("A112").Value = "Kleuren"
$ws.Range("B112").Value = 1
$ws.Range("C112").Value = "grijs"
$ws.Range("D112").Value = "gris"
$ws.Range("B113").Value = 1
$ws.Range("D113").Value = "mono"
$ws.Range("C113").Value = "aap"
$ws.Range("B114").Value = 1
$ws.Range("D114").Value = "oso"
$ws.Range("C114").Value = "beer"
$ws.Range("B115").Value = 1
$ws.Range("C115").Value = "olifant"
$ws.Range("B116").Value = 1
$ws.Range("C116").Value = "kangoeroe"
$ws.Range("B117").Value = 1
$ws.Range("C117").Value = "giraffe"
$ws.Range("B118").Value = 1
$ws.Range("C118").Value = "neushoorn"
$ws.Range("B119").Value = 1
$ws.Range("C119").Value = "nijlpaard"
$ws.Range("B120").Value = 1
$ws.Range("C120").Value = "leeuw"
$ws.Range("B121").Value = 1
$ws.Range("C121").Value = "tijger"
$ws.Range("B122").Value = 1
$ws.Range("C122").Value = "kameel"
$ws.Range("B123").Value = 1
$ws.Range("C123").Value = "hert"
$ws.Range("B124").Value = 1
$ws.Range("C124").Value = "slang"
$ws.Range("B125").Value = 1
$ws.Range("C125").Value = "krokodil"
$ws.Range("B126").Value = 1
$ws.Range("C126").Value = "schildpad"
$ws.Range("B127").Value = 1
$ws.Range("C127").Value = "hagedis"
$ws.Range("B128").Value = 1
$ws.Range("C128").Value = "salamander"
$ws.Range("B129").Value = 1
$ws.Range("C129").Value = "kikker"
$ws.Range("B130").Value = 1
$ws.Range("C130").Value = "pad"
$ws.Range("B131").Value = 1
$ws.Range("C131").Value = "eend"
$ws.Range("B132").Value = 1
$ws.Range("C132").Value = "zwaan"
$ws.Range("B133").Value = 1
$ws.Range("C133").Value = "gans"
$ws.Range("B134").Value = 1
$ws.Range("C134").Value = "mus"
$ws.Range("B135").Value = 1
$ws.Range("C135").Value = "duif"
$ws.Range("B136").Value = 1
$ws.Range("C136").Value = "merel"
$ws.Range("B137").Value = 1
$ws.Range("C137").Value = "meeuw"
$ws.Range("B138").Value = 1
$ws.Range("C138").Value = "roofvogel"
$ws.Range("B139").Value = 1
$ws.Range("C139").Value = "uil"
$ws.Range("B140").Value = 1
$ws.Range("C140").Value = "struisvogel"
$ws.Range("B141").Value = 1
$ws.Range("C141").Value = "pinguin"
$ws.Range("B142").Value = 1
$ws.Range("C142").Value = "papegaai"
$ws.Range("B143").Value = 1
$ws.Range("C143").Value = "zeehond"
$ws.Range("B144").Value = 1
$ws.Range("C144").Value = "dolfijn"
$ws.Range("B145").Value = 1
$ws.Range("C145").Value = "haai"
$ws.Range("B146").Value = 1
$ws.Range("C146").Value = "vis"
$ws.Range("B147").Value = 1
$ws.Range("C147").Value = "orka"
$ws.Range("B148").Value = 1
$ws.Range("C148").Value = "spin"
$ws.Range("B149").Value = 1
$ws.Range("C149").Value = "bij"
$ws.Range("B150").Value = 1
$ws.Range("C150").Value = "wesp"
$ws.Range("B151").Value = 1
$ws.Range("C151").Value = "mug"
$ws.Range("B152").Value = 1
$ws.Range("C152").Value = "vlieg"
$ws.Range("B153").Value = 1
$ws.Range("C153").Value = "mier"
$ws.Range("B154").Value = 1
$ws.Range("C154").Value = "rups"
$ws.Range("B155").Value = 1
$ws.Range("C155").Value = "vliender"
$ws.Range("B156").Value = 1
$ws.Range("C156").Value = "kever"
$ws.Range("B157").Value = 1
$ws.Range("C157").Value = "hond"
$ws.Range("B158").Value = 1
$ws.Range("C158").Value = "kat"
$ws.Range("B159").Value = 1
$ws.Range("C159").Value = "konijn"
$ws.Range("B160").Value = 1
$ws.Range("C160").Value = "cavia"
$ws.Range("B161").Value = 1
$ws.Range("C161").Value = "muis"
$ws.Range("B162").Value = 1
$ws.Range("C162").Value = "rat"
$ws.Range("B163").Value = 1
$ws.Range("C163").Value = "hamster"
$ws.Range("B164").Value = 1
$ws.Range("C164").Value = "goudvis"
$ws.Range("B165").Value = 1
$ws.Range("C165").Value = "koe"
$ws.Range("B166").Value = 1
$ws.Range("C166").Value = "stier"
$ws.Range("B167").Value = 1
$ws.Range("C167").Value = "paard"
$ws.Range("B168").Value = 1
$ws.Range("C168").Value = "varken"
$ws.Range("B169").Value = 1
$ws.Range("C169").Value = "ezel"
$ws.Range("B170").Value = 1
$ws.Range("C170").Value = "schaap"
$ws.Range("B171").Value = 1
$ws.Range("C171").Value = "kip"
$ws.Range("B172").Value = 1
$ws.Range("C172").Value = "haan"
$ws.Range("D115").Value = "elefante"
$ws.Range("D116").Value = "canguro"
$ws.Range("D117").Value = "jirafa"
$ws.Range("D118").Value = "rinoceronte"
$ws.Range("D119").Value = "hipopótamo"
$ws.Range("D120").Value = "león"
$ws.Range("D121").Value = "tigre"
$ws.Range("D122").Value = "camello"
$ws.Range("D123").Value = "venado"
$ws.Range("D124").Value = "serpiente"
$ws.Range("D125").Value = "cocodrilo"
$ws.Range("D126").Value = "tortuga"
$ws.Range("D127").Value = "lagargo"
$ws.Range("D128").Value = "salamandra"
$ws.Range("D129").Value = "rana"
$ws.Range("D130").Value = "sapo"
$ws.Range("D131").Value = "pato"
$ws.Range("D132").Value = "cisne"
$ws.Range("D133").Value = "ganso"
$ws.Range("D134").Value = "gorrión"
$ws.Range("D135").Value = "paloma"
$ws.Range("D136").Value = "mirlo"
$ws.Range("D137").Value = "gaviota"
$ws.Range("D138").Value = "ave de presa"
$ws.Range("D139").Value = "búho"
$ws.Range("D140").Value = "avestruz"
$ws.Range("D141").Value = "pingüino"
$ws.Range("D142").Value = "loro"
$ws.Range("D143").Value = "foca"
$ws.Range("D144").Value = "delfín"
$ws.Range("D145").Value = "tiburón"
$ws.Range("D146").Value = "pez"
$ws.Range("D147").Value = "orca"
$ws.Range("D148").Value = "araña"
$ws.Range("D149").Value = "abeja"
$ws.Range("D150").Value = "avispa"
$ws.Range("D151").Value = "mosquito"
$ws.Range("D152").Value = "mosca"
$ws.Range("D153").Value = "hormiga"
$ws.Range("D154").Value = "oruga"
$ws.Range("D155").Value = "mariposa"
$ws.Range("D156").Value = "escarabajo"
$ws.Range("D157").Value = "perro"
$ws.Range("D158").Value = "gato"
$ws.Range("D159").Value = "conejo"
$ws.Range("D160").Value = "conejillo de indias"
$ws.Range("D161").Value = "ratón"
$ws.Range("D162").Value = "rata"
$ws.Range("D163").Value = "hamster"
$ws.Range("D164").Value = "pez dorado"
$ws.Range("D165").Value = "vaca"
$ws.Range("D166").Value = "toro"
$ws.Range("D167").Value = "caballo"
$ws.Range("D168").Value = "cerdo"
$ws.Range("D169").Value = "burro"
$ws.Range("D170").Value = "oveja"
$ws.Range("D171").Value = "pollo"
$ws.Range("D172").Value = "gallo"
$ws.Range("A113").Value = "Dieren"
$ws.Range("A114").Value = "Dieren"
$ws.Range("A115").Value = "Dieren"
$ws.Range("A116").Value = "Dieren"
$ws.Range("A117").Value = "Dieren"
$ws.Range("A118").Value = "Dieren"
$ws.Range("A119").Value = "Dieren"
$ws.Range("A120").Value = "Dieren"
$ws.Range("A121").Value = "Dieren"
$ws.Range("A122").Value = "Dieren"
$ws.Range("A123").Value = "Dieren"
$ws.Range("A124").Value = "Dieren"
$ws.Range("A125").Value = "Dieren"
$ws.Range("A126").Value = "Dieren"
$ws.Range("A127").Value = "Dieren"
$ws.Range("A128").Value = "Dieren"
$ws.Range("A129").Value = "Dieren"
$ws.Range("A130").Value = "Dieren"
$ws.Range("A131").Value = "Dieren"
$ws.Range("A132").Value = "Dieren"
$ws.Range("A133").Value = "Dieren"
$ws.Range("A134").Value = "Dieren"
$ws.Range("A135").Value = "Dieren"
$ws.Range("A136").Value = "Dieren"
$ws.Range("A137").Value = "Dieren"
$ws.Range("A138").Value = "Dieren"
$ws.Range("A139").Value = "Dieren"
$ws.Range("A140").Value = "Dieren"
$ws.Range("A141").Value = "Dieren"
$ws.Range("A142").Value = "Dieren"
$ws.Range("A143").Value = "Dieren"
$ws.Range("A144").Value = "Dieren"
$ws.Range("A145").Value = "Dieren"
$ws.Range("A146").Value = "Dieren"
$ws.Range("A147").Value = "Dieren"
$ws.Range("A148").Value = "Dieren"
$ws.Range("A149").Value = "Dieren"
$ws.Range("A150").Value = "Dieren"
$ws.Range("A151").Value = "Dieren"
$ws.Range("A152").Value = "Dieren"
$ws.Range("A153").Value = "Dieren"
$ws.Range("A154").Value = "Dieren"
$ws.Range("A155").Value = "Dieren"
$ws.Range("A156").Value = "Dieren"
$ws.Range("A157").Value = "Dieren"
$ws.Range("A158").Value = "Dieren"
$ws.Range("A159").Value = "Dieren"
$ws.Range("A160").Value = "Dieren"
$ws.Range("A161").Value = "Dieren"
$ws.Range("A162").Value = "Dieren"
$ws.Range("A163").Value = "Dieren"
$ws.Range("A164").Value = "Dieren"
$ws.Range("A165").Value = "Dieren"
$ws.Range("A166").Value = "Dieren"
$ws.Range("A167").Value = "Dieren"
$ws.Range("A168").Value = "Dieren"
$ws.Range("A169").Value = "Dieren"
$ws.Range("A170").Value = "Dieren"
$ws.Range("A171").Value = "Dieren"
$ws.Range("A172").Value = "Dieren"

# Update the active selection to mirror the final edit session
$ws.Activate()
$ws.Range("B113:B172").Select()

Write-Output "done"